$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Đơn sale chính"
# Insert a new row 3 (pushing the old "Tổng" row down to row 4) and fill it
# with the new HD-LUXURY/598 order, then update the totals row.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = "HD-LUXURY"
$ws1.Range("B3").Value = 598
$ws1.Range("C3").NumberFormat = "@"
$ws1.Range("C3").Value = "07-28-2024"
$ws1.Range("D3").Value = "CẦN THƠ"
$ws1.Range("E3").Value = "Thị Minh"
$ws1.Range("F3").Value = "Cá nhân"
$ws1.Range("G3").Value = "Cọc Dịch Vụ Làm Đẹp"
$ws1.Range("H3").Value = 200000
$ws1.Range("K3").Value = 200000
$ws1.Range("L3").Value = 200000
$ws1.Range("M3").Value = 0
$ws1.Range("N3").Value = 0

# Update the totals row (now row 4)
$ws1.Range("B4").Value = 2
$ws1.Range("H4").Value = 2000000
$ws1.Range("K4").Value = 2000000
$ws1.Range("L4").Value = 2000000

# ---------------------------------------------------------------------------
# Sheet 2: "Lương"
# Insert a new row 2 (pushing everything down by one) so the payroll
# breakdown gains a "Tổng công / Phụ cấp / Lương công tác tại CẦN THƠ" block
# at the top, then relabel/reassign the shifted rows and append the new
# trailing rows for SÓC TRĂNG details and grand totals.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(2).Insert()

$ws2.Range("A2").Value = "Tổng công tại CẦN THƠ"
$ws2.Range("B2").Value = 0
$ws2.Range("A3").Value = "Phụ cấp tại CẦN THƠ"
$ws2.Range("B3").Value = 0
$ws2.Range("A4").Value = "Lương công tác tại CẦN THƠ"
$ws2.Range("B4").Value = 0
$ws2.Range("B5").ClearContents()
$ws2.Range("A13").Value = "Tổng công tại LONG XUYÊN"
$ws2.Range("B13").Value = 28
$ws2.Range("A14").Value = "Phụ cấp tại LONG XUYÊN"
$ws2.Range("B14").Value = 980000
$ws2.Range("A15").Value = "Lương cơ bản tại LONG XUYÊN"
$ws2.Range("B15").Value = 4000000
$ws2.Range("A16").Value = "Chiết khấu sale chính tại LONG XUYÊN"
$ws2.Range("B16").Value = 180000
$ws2.Range("A17").Value = "Chiết khấu sale phụ tại LONG XUYÊN"
$ws2.Range("A18").Value = "Đơn 1 bác sĩ tại LONG XUYÊN"
$ws2.Range("A19").Value = "Đơn 2 bác sĩ tại LONG XUYÊN"
$ws2.Range("A20").Value = "Công phụ phẫu 1 tại LONG XUYÊN"
$ws2.Range("A21").Value = "Công phụ phẫu 2 tại LONG XUYÊN"
$ws2.Range("B21").Value = 0
$ws2.Range("A22").Value = "Ứng lương tại LONG XUYÊN"
$ws2.Range("A23").Value = "Tổng công tại SÓC TRĂNG"
$ws2.Range("A24").Value = "Phụ cấp tại SÓC TRĂNG"
$ws2.Range("A25").Value = "Lương công tác tại SÓC TRĂNG"
$ws2.Range("A26").Value = "Lương cơ bản tại SÓC TRĂNG"
$ws2.Range("B26").ClearContents()
$ws2.Range("A27").Value = "Chiết khấu sale chính tại SÓC TRĂNG"
$ws2.Range("A28").Value = "Chiết khấu sale phụ tại SÓC TRĂNG"
$ws2.Range("A29").Value = "Đơn 1 bác sĩ tại SÓC TRĂNG"
$ws2.Range("A30").Value = "Đơn 2 bác sĩ tại SÓC TRĂNG"
$ws2.Range("B30").Value = 0
$ws2.Range("A31").Value = "Công phụ phẫu 1 tại SÓC TRĂNG"
$ws2.Range("A32").Value = "Công phụ phẫu 2 tại SÓC TRĂNG"
$ws2.Range("B32").Value = 0
$ws2.Range("A33").Value = "Ứng lương tại SÓC TRĂNG"
$ws2.Range("B33").Value = 0
$ws2.Range("A34").Value = "Tổng lương tại CẦN THƠ"
$ws2.Range("B34").Value = 0
$ws2.Range("A35").Value = "Tổng lương tại LONG XUYÊN"
$ws2.Range("B35").Value = 5160000
$ws2.Range("A36").Value = "Tổng lương tại SÓC TRĂNG"
$ws2.Range("B36").Value = 0
$ws2.Range("A37").Value = "Tổng lương"
$ws2.Range("B37").Value = 5160000
